# Add a new slide (13) at the end of the deck, using the same
# "Title and Content" layout (ppLayoutObject == 16) as every other
# content slide in this deck, and fill in its placeholders.

$p = $ppt.ActivePresentation

$s = $p.Slides.Add($p.Slides.Count + 1, 16)

# Title placeholder.
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Download complete presentation here"

# Body / content placeholder with the repo link, sized to 40pt like the
# author did (this also turns on normAutofit, matching the source deck).
$body = $s.Shapes.Item(2)
$body.TextFrame.TextRange.Text = "https://github.com/sunnyhlopez/cms_class"
$body.TextFrame.TextRange.Font.Size = 40

# Bring in the (empty) footer placeholder the same way every other
# slide in this deck has one: copy the existing footer placeholder
# shape from the last slide and paste it onto the new slide, rather
# than re-creating it from scratch.
$srcFooter = $p.Slides.Item($p.Slides.Count - 1).Shapes.Item(3)
$srcFooter.Copy()
$s.Shapes.Paste() | Out-Null
